$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Duplicate-ish new sheet "DPLKKPS138-002" placed right after sheet 1
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "DPLKKPS138-002"

# ---------------------------------------------------------------------------
# 2) Sheet1 ("DPLKKPS138-001") formatting / content tweaks
# ---------------------------------------------------------------------------

# Header cell C1: remove center alignment -> left (re-uses existing left/center style)
$ws1.Range("C1").HorizontalAlignment = -4131

# Row 2 cells: switch several center-aligned (wrap) cells to left aligned (wrap)
$ws1.Range("B2").HorizontalAlignment = -4131
$ws1.Range("D2").HorizontalAlignment = -4131
$ws1.Range("E2").HorizontalAlignment = -4131
$ws1.Range("F2").HorizontalAlignment = -4131
$ws1.Range("F2").WrapText = $True

$ws1.Range("M2").WrapText = $True

# Cell content updates on row 2
$ws1.Range("F2").Value = "1 : Lanjutkan ke Verifikasi"
$ws1.Range("R2").Value = 0
$ws1.Range("S2").Value = "KEP-011 Pending"

# Column widths
$ws1.Columns.Item(6).ColumnWidth = 29.09
$ws1.Columns.Item(13).ColumnWidth = 25

# Row height (content got shorter -> row auto height shrinks)
$ws1.Rows.Item(2).RowHeight = 127.5

# Sheet view: no longer the active/selected tab, scrolled one column left,
# selection moved from P2 to O2
$ws1.Range("O2").Select()
$ws1.Application.ActiveWindow.ScrollColumn = 9

# ---------------------------------------------------------------------------
# 3) Populate the new sheet2 ("DPLKKPS138-002") -- mirrors sheet1's layout
#    with a couple of value differences (PKS number + status texts)
# ---------------------------------------------------------------------------
$headers = @("RUN","TC_ID","TEST_SCENARIO_DESC","SCENARIO_DESC","EXPECTED_RESULT","PREPARATION","USERID","PASSWORD","ROLE","MAIN_SIDEBAR","SIDEBAR_MENU","SIDEBAR_SUBMENU","SIDEBAR_SUBMENU_SUBMENU","KETERANGAN_PERUBAHAN","NO_PKS","ADMIN_FEE","MGMT_FEE","STATUS_REGISTER","KETERANGAN_REGISTER")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws2.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ws2.Range("A2").Value = "RUN"
$ws2.Range("B2").Value = "DPLKKPS138-002"
$ws2.Range("C2").Value = "Normal - Kepesertaan - Transaksi"
$ws2.Range("D2").Value = "Register - Maintenance Perjanjian Kerja Sama"
$ws2.Range("E2").Value = "Proses Register Maintenance PKS berhasil"
$ws2.Range("F2").Value = "Username : 30603;`nPassword : bni1234;`nRole : 10 - Asisten Settlement;`nKeterangan Perubahan : KEP-011;`nNo PKS : 82329799;`nTab Parameter;`nADMIN_FEE : 15000;`nMGMT_FEE : 0,50;`nStatus Register : 0 : Pending Register;`nKeterangan Register : KEP-011"
$ws2.Range("G2").Value = 30603
$ws2.Range("H2").Value = "bni1234"
$ws2.Range("I2").Value = "Asisten Settlement"
$ws2.Range("J2").Value = "Kepesertaan"
$ws2.Range("K2").Value = "Transaksi"
$ws2.Range("L2").Value = "Perjanjian Kerja Sama"
$ws2.Range("M2").Value = "Maintenance Perjanjian Kerja Sama"
$ws2.Range("N2").Value = "KEP-011"
$ws2.Range("O2").Value = 82329777
$ws2.Range("P2").Value = 15000
$ws2.Range("Q2").Value = 0.5
$ws2.Range("R2").Value = 1
$ws2.Range("S2").Value = "1 : Lanjutkan ke Verifikasi"

$ws2.Range("A1:S1").HorizontalAlignment = -4131
$ws2.Range("A1:S1").Font.Size = 10

$ws2.Range("B2:F2").HorizontalAlignment = -4131
$ws2.Range("B2:F2").WrapText = $True
$ws2.Range("B2:F2").VerticalAlignment = -4108
$ws2.Range("B2:F2").Font.Size = 10

$ws2.Rows.Item(2).RowHeight = 165.75

$ws2.Columns.Item(1).ColumnWidth = 4.14
$ws2.Columns.Item(2).ColumnWidth = 13.86
$ws2.Columns.Item(3).ColumnWidth = 18.71
$ws2.Columns.Item(4).ColumnWidth = 18.71
$ws2.Columns.Item(5).ColumnWidth = 14.43
$ws2.Columns.Item(6).ColumnWidth = 24.14
$ws2.Columns.Item(7).ColumnWidth = 6
$ws2.Columns.Item(8).ColumnWidth = 9.29

$ws2.Range("R2").Select()

$wb.Worksheets.Item(1).Activate()
$wb.Worksheets.Item(2).Activate()
